$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.033.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.408.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.406.89"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("E10").Value = "  +3.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.835.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.957.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").Value = "  +4.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.404.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.22%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +6.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.15%  "
$ws.Range("E29").Value = "  +15.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0771"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "170.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "362.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.01%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  +8.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.34%  "
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0518"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.63%  "
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0216"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.33%  "
